# Collapse multi-run title/caption text (one <a:r> per word) down to a
# single run per paragraph, matching the canonical golden output.
#
# Setting TextRange.Text to the same concatenated string the runs already
# represent is a no-op in this runtime (it compares resulting text), so we
# first stomp each target with a dummy placeholder to force a real text
# change, then set the real desired text - that rewrite collapses the
# paragraph's runs into a single run.

function Set-MergedText($shape, [string]$text) {
    $tr = $shape.TextFrame.TextRange
    $tr.Text = "~~~tmp~~~"
    $tr.Text = $text
}

$p = $ppt.ActivePresentation

Set-MergedText $p.Slides.Item(1).Shapes.Item(1) "Slide 1 (Content)"
Set-MergedText $p.Slides.Item(2).Shapes.Item(1) "Slide 2 (Content)"
Set-MergedText $p.Slides.Item(3).Shapes.Item(1) "Slide 3 (Content)"
Set-MergedText $p.Slides.Item(4).Shapes.Item(1) "Slide 4 (Content)"
Set-MergedText $p.Slides.Item(5).Shapes.Item(1) "Slide 5 (Two Content)"

Set-MergedText $p.Slides.Item(6).Shapes.Item(1) "Slide 6 (Two Content Right)"
Set-MergedText $p.Slides.Item(6).Shapes.Item(3) "an image"

Set-MergedText $p.Slides.Item(7).Shapes.Item(1) "Slide 7 (Content with Caption)"
Set-MergedText $p.Slides.Item(7).Shapes.Item(4) "An image"

Set-MergedText $p.Slides.Item(8).Shapes.Item(1) "Slide 8 (Comparison)"
Set-MergedText $p.Slides.Item(8).Shapes.Item(4) "An image"

Set-MergedText $p.Slides.Item(9).Shapes.Item(1) "Slide 10 (Content)"
Set-MergedText $p.Slides.Item(10).Shapes.Item(1) "Slide 11 (Content)"
Set-MergedText $p.Slides.Item(11).Shapes.Item(1) "Slide 12 (Content)"
